$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.856.24"
$ws.Range("E2").Value = "  -1.52%  "

$ws.Range("D3").Value = "3.170.43"

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.53"
$ws.Range("E5").Value = "  -2.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.97"
$ws.Range("E6").Value = "  -6.16%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.167.49"
$ws.Range("E8").Value = "  -4.53%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("E9").Value = "  -0.87%  "

$ws.Range("E10").Value = "  -6.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.22"
$ws.Range("E11").Value = "  -5.95%  "

$ws.Range("E12").Value = "  -4.05%  "

$ws.Range("E13").Value = "  -5.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.61"
$ws.Range("E14").Value = "  -1.32%  "

$ws.Range("D15").Value = "3.692.09"
$ws.Range("E15").Value = "  -4.51%  "

$ws.Range("E16").Value = "  -1.24%  "

$ws.Range("D17").Value = "3.170.11"
$ws.Range("E17").Value = "  -4.53%  "

$ws.Range("D18").Value = "62.838.71"
$ws.Range("E18").Value = "  -1.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.53"
$ws.Range("E19").Value = "  -4.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "459.24"
$ws.Range("E20").Value = "  -4.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.98"
$ws.Range("E21").Value = "  -0.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.691"
$ws.Range("E22").Value = "  -6.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.60"
$ws.Range("E23").Value = "  -4.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.31"
$ws.Range("E24").Value = "  -4.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.61"
$ws.Range("E25").Value = "  -2.86%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("E28").Value = "  -4.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.67"
$ws.Range("E29").Value = "  -7.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.74"
$ws.Range("E30").Value = "  -5.92%  "

$ws.Range("E31").Value = "  -6.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.07"
$ws.Range("E32").Value = "  -6.30%  "

$ws.Range("E33").Value = "  -4.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.35"
$ws.Range("E34").Value = "  -7.06%  "

$ws.Range("E35").Value = "  -6.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.79"
$ws.Range("E36").Value = "  -5.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.27"
$ws.Range("E37").Value = "  -2.26%  "

$ws.Range("D38").Value = "0.0₃0701"
$ws.Range("E38").Value = "  -6.02%  "

$ws.Range("E39").Value = "  -3.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "403.87"
$ws.Range("E40").Value = "  -7.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.04"
$ws.Range("E41").Value = "  -3.91%  "

$ws.Range("E42").Value = "  -5.57%  "

$ws.Range("E43").Value = "  -5.29%  "

$ws.Range("D44").Value = "2.790.06"
$ws.Range("E44").Value = "  -10.82%  "

$ws.Range("E45").Value = "  -6.61%  "

$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.10"
$ws.Range("E47").Value = "  -6.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.69"
$ws.Range("E48").Value = "  +0.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.12"
$ws.Range("E49").Value = "  -4.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.32"
$ws.Range("E50").Value = "  -6.55%  "

$ws.Range("E51").Value = "  -2.34%  "
